$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 381, shifting existing rows 381-409 down to 382-410
$ws.Rows("381").Insert()

# Populate the newly inserted row 381 with the new weekly record
$ws.Range("A381").Value = 3
$ws.Range("B381").Value = "Femacal de La Calera"
$ws.Range("C381").Value = "Coquimbo"
$ws.Range("D381").Value = 44714
$ws.Range("E381").Value = 5
$ws.Range("F381").Value = "Fruta"
$ws.Range("G381").Value = 100108
$ws.Range("H381").Value = "Tropicales y subtropicales"
$ws.Range("I381").Value = 100108002
$ws.Range("J381").Value = "Mango"
$ws.Range("K381").Value = "Sin especificar"
$ws.Range("L381").Value = "Primera"
$ws.Range("M381").Value = 456
$ws.Range("N381").Value = 9000
$ws.Range("O381").Value = 9000
$ws.Range("P381").Value = 9000
$ws.Range("Q381").Value = "$/bandeja 4 kilos"
$ws.Range("R381").Value = "Brasil"
$ws.Range("S381").Value = 2250
$ws.Range("T381").Value = 4
